$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.928.55"
$ws.Range("E2").Value = "  +4.15%  "
$ws.Range("D3").Value = "3.249.69"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'396.03"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "'108.68"
$ws.Range("E6").Value = "  -1.20%  "
$ws.Range("D7").Value = "'0.589"
$ws.Range("E7").Value = "  +7.11%  "
$ws.Range("D8").Value = "3.248.08"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("D10").Value = "'0.626"
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("D11").Value = "'39.25"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "'0.0993"
$ws.Range("E12").Value = "  +11.82%  "
$ws.Range("E13").Value = "  +2.16%  "
$ws.Range("D14").Value = "3.773.17"
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("D15").Value = "'8.35"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").Value = "'19.15"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "3.245.68"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("E18").Value = "  -3.13%  "
$ws.Range("D19").Value = "'10.75"
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "56.816.66"
$ws.Range("E20").Value = "  +4.22%  "
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  +8.18%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").Value = "'295.28"
$ws.Range("E24").Value = "  +6.91%  "
$ws.Range("D25").Value = "'74.49"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").Value = "'28.14"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "'4.36"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").Value = "'7.66"
$ws.Range("E29").Value = "  -4.60%  "
$ws.Range("D30").Value = "'7.24"
$ws.Range("E30").Value = "  -3.89%  "
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'11.25"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("E34").Value = "  -2.93%  "
$ws.Range("D35").Value = "'39.87"
$ws.Range("E35").Value = "  +8.69%  "
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "'51.58"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -3.88%  "
$ws.Range("D41").Value = "'2.95"
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'139.25"
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("D43").Value = "'0.122"
$ws.Range("E43").Value = "  +3.88%  "
$ws.Range("E44").Value = "  -2.18%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'17.14"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'3.97"
$ws.Range("E46").Value = "  -2.81%  "
$ws.Range("D47").Value = "'0.281"
$ws.Range("E47").Value = "  -3.96%  "
$ws.Range("D48").Value = "'22.28"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "'2.16"
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").Value = "2.165.21"
$ws.Range("E50").Value = "  +3.37%  "
$ws.Range("E51").Value = "  -6.54%  "
